$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: date moves back one day (was "12월 9일", now "12월 8일").
# The "Activity" text in F34 stays the same ("최종결과보고서 작성").
$ws.Range("A34").Value = "12월 8일"
$ws.Range("F34").Value = "최종결과보고서 작성"

# Row 35 was blank; it now holds the entry that used to live in row 34
# (date "12월 9일") together with its own time-log values.
$ws.Range("A35").Value = "12월 9일"
$ws.Range("B35").Value = 0.041666666666666664
$ws.Range("C35").Value = 0.29166666666666669
$ws.Range("D35").Value = 60
$ws.Range("E35").Value = 360
$ws.Range("F35").Value = "최종결과보고서 작성"

# F35 picks up the same font used by the other "Activity" cells
# (style index 23 in the original file: Arial Unicode MS, size 10).
$ws.Range("F35").Font.Name = "Arial Unicode MS"
$ws.Range("F35").Font.Size = 10

# Update the saved view state: the selection now spans the
# newly-completed two-row block.
$ws.Range("A34:F35").Select()
